# Applies the cryptos-list refresh described in the commit
# "Updated cryptos list on Sat Oct 28 18:48:37 UTC 2023 with GitHub Actions".
# All target values are written as text (leading "'" forces text so that
# values such as "225.98" or "0.0520" are not reinterpreted as numbers,
# matching the original inline-string cell contents).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.179.71"
$ws.Range("E2").Value = "'  +1.65%  "
$ws.Range("D3").Value = "'1.785.17"
$ws.Range("E3").Value = "'  +1.02%  "
$ws.Range("E4").Value = "'  +0.09%  "
$ws.Range("D5").Value = "'225.98"
$ws.Range("E5").Value = "'  +1.22%  "
$ws.Range("E6").Value = "'  +0.60%  "
$ws.Range("E7").Value = "'  +0.17%  "
$ws.Range("D8").Value = "'31.79"
$ws.Range("E8").Value = "'  +0.29%  "
$ws.Range("D9").Value = "'0.292"
$ws.Range("E9").Value = "'  +1.45%  "
$ws.Range("E10").Value = "'  +0.44%  "
$ws.Range("E11").Value = "'  +1.20%  "
$ws.Range("D12").Value = "'2.041.90"
$ws.Range("E12").Value = "'  +1.08%  "
$ws.Range("D13").Value = "'11.02"
$ws.Range("E13").Value = "'  -0.81%  "
$ws.Range("D14").Value = "'1.782.52"
$ws.Range("B15").Value = "'Polygon"
$ws.Range("C15").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.623"
$ws.Range("E15").Value = "'  +2.51%  "
$ws.Range("B16").Value = "'WrappedBTC"
$ws.Range("C16").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "'34.122.69"
$ws.Range("E16").Value = "'  +1.34%  "
$ws.Range("E17").Value = "'  +1.74%  "
$ws.Range("E18").Value = "'  +2.75%  "
$ws.Range("D19").Value = "'246.42"
$ws.Range("E19").Value = "'  +4.07%  "
$ws.Range("D20").Value = "'0.0₃0778"
$ws.Range("E20").Value = "'  +0.94%  "
$ws.Range("E21").Value = "'  -0.06%  "
$ws.Range("D22").Value = "'10.93"
$ws.Range("E22").Value = "'  +3.73%  "
$ws.Range("E23").Value = "'  +2.48%  "
$ws.Range("E24").Value = "'  +0.31%  "
$ws.Range("D25").Value = "'161.49"
$ws.Range("E25").Value = "'  +1.48%  "
$ws.Range("D26").Value = "'7.19"
$ws.Range("E26").Value = "'  +2.74%  "
$ws.Range("D27").Value = "'16.33"
$ws.Range("E27").Value = "'  +1.59%  "
$ws.Range("D28").Value = "'0.114"
$ws.Range("E28").Value = "'  +2.09%  "
$ws.Range("E29").Value = "'  +0.33%  "
$ws.Range("E30").Value = "'  +1.11%  "
$ws.Range("D31").Value = "'0.0520"
$ws.Range("E31").Value = "'  +1.83%  "
$ws.Range("D32").Value = "'3.69"
$ws.Range("E32").Value = "'  +3.00%  "
$ws.Range("D33").Value = "'3.63"
$ws.Range("E33").Value = "'  +4.16%  "
$ws.Range("D34").Value = "'1.80"
$ws.Range("E34").Value = "'  +1.35%  "
$ws.Range("D35").Value = "'1.449.95"
$ws.Range("E35").Value = "'  +5.26%  "
$ws.Range("D36").Value = "'0.654"
$ws.Range("E36").Value = "'  +1.25%  "
$ws.Range("E37").Value = "'  +8.81%  "
$ws.Range("E38").Value = "'  +4.23%  "
$ws.Range("E39").Value = "'  +1.68%  "
$ws.Range("D40").Value = "'80.31"
$ws.Range("E40").Value = "'  +3.42%  "
$ws.Range("D41").Value = "'2.37"
$ws.Range("E41").Value = "'  +0.72%  "
$ws.Range("D42").Value = "'0.923"
$ws.Range("E42").Value = "'  +2.47%  "
$ws.Range("D43").Value = "'2.69"
$ws.Range("E43").Value = "'  +0.90%  "
$ws.Range("D44").Value = "'13.51"
$ws.Range("E44").Value = "'  +1.07%  "
$ws.Range("E45").Value = "'  +4.85%  "
$ws.Range("E46").Value = "'  +2.28%  "
$ws.Range("E47").Value = "'  -0.53%  "
$ws.Range("D48").Value = "'0.0₆0136"
$ws.Range("E48").Value = "'  -0.90%  "
$ws.Range("D49").Value = "'1.944.25"
$ws.Range("E49").Value = "'  +1.61%  "
$ws.Range("D50").Value = "'105.89"
$ws.Range("E50").Value = "'  -0.74%  "
$ws.Range("E51").Value = "'  +0.11%  "
